$d = $word.ActiveDocument

# Title: "Play Natural Powers for Free - Slot Game Review" -> "Play Natural Powers for Free 2021"
# (occurs twice: the H1 heading near the top, and the bold "recap" line near the
# bottom - Replace:=wdReplaceAll (2) takes care of every occurrence in one call)
$d.Content.Find.Execute("Play Natural Powers for Free - Slot Game Review", $true, $false, $false, $false, $false, $true, 1, $false, "Play Natural Powers for Free 2021", 2)

# "What we like" bullet list
$d.Content.Find.Execute("Rolling reels format keeps gameplay exciting", $true, $false, $false, $false, $false, $true, 1, $false, "Unique rolling reels format adds excitement to gameplay", 2)
$d.Content.Find.Execute("Expanding symbols add a unique twist", $true, $false, $false, $false, $false, $true, 1, $false, "Appealing superheroine theme for fans of superhero stories", 2)
$d.Content.Find.Execute("Free spin/bonus round options are varied", $true, $false, $false, $false, $false, $true, 1, $false, "Stunning graphics and sound effects enhance the overall experience", 2)
$d.Content.Find.Execute("Stunning graphics enhance the game's atmosphere", $true, $false, $false, $false, $false, $true, 1, $false, "Free spin/bonus rounds offer a variety of choices and multipliers", 2)

# "What we don't like" bullet list
$d.Content.Find.Execute("No progressive jackpot feature", $true, $false, $false, $false, $false, $true, 1, $false, "Limited fixed paylines may not appeal to players looking for more flexibility", 2)
$d.Content.Find.Execute("No multipliers during base gameplay", $true, $false, $false, $false, $false, $true, 1, $false, "Lack of variety in symbol design besides the superheroines", 2)

# Closing italic summary line
$d.Content.Find.Execute("Experience the power of Natural Powers. Read our slot game review and play for free. Engage in the thrilling gameplay mechanics with four expanding symbols.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of the online slot game Natural Powers and play it for free today!", 2)
